# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Rebuilds the worker/period detail table (rows 16-33) on Hoja1.
# Each worker now carries periods 2103-2109 (Marjorie only 2103-2106),
# with refreshed "Valor Mora" (F) / "Salario Basico" (G) amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# MARJORIE VARELA TORRES (CC 45692124) - rows 16-19
$ws.Range("C16").Value = "45692124"
$ws.Range("D16").Value = "MARJORIE VARELA TORRES"
$ws.Range("E16").Value = "2106"
$ws.Range("F16").Value = 48000
$ws.Range("G16").Value = 1700000

$ws.Range("C17").Value = "45692124"
$ws.Range("D17").Value = "MARJORIE VARELA TORRES"
$ws.Range("E17").Value = "2105"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1700000

$ws.Range("C18").Value = "45692124"
$ws.Range("D18").Value = "MARJORIE VARELA TORRES"
$ws.Range("E18").Value = "2104"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1700000

$ws.Range("C19").Value = "45692124"
$ws.Range("D19").Value = "MARJORIE VARELA TORRES"
$ws.Range("E19").Value = "2103"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1700000

# WINYBEL SANDRY DE LA HOZ MUÑOZ (CC 1051890181) - rows 20-26
$ws.Range("C20").Value = "1051890181"
$ws.Range("D20").Value = "WINYBEL SANDRY DE LA HOZ MUÑOZ"
$ws.Range("E20").Value = "2109"
$ws.Range("F20").Value = 24227
$ws.Range("G20").Value = 908526

$ws.Range("C21").Value = "1051890181"
$ws.Range("D21").Value = "WINYBEL SANDRY DE LA HOZ MUÑOZ"
$ws.Range("E21").Value = "2108"
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 908526

$ws.Range("C22").Value = "1051890181"
$ws.Range("D22").Value = "WINYBEL SANDRY DE LA HOZ MUÑOZ"
$ws.Range("E22").Value = "2107"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 908526

$ws.Range("C23").Value = "1051890181"
$ws.Range("D23").Value = "WINYBEL SANDRY DE LA HOZ MUÑOZ"
$ws.Range("E23").Value = "2106"
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 908526

$ws.Range("C24").Value = "1051890181"
$ws.Range("D24").Value = "WINYBEL SANDRY DE LA HOZ MUÑOZ"
$ws.Range("E24").Value = "2105"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 908526

$ws.Range("C25").Value = "1051890181"
$ws.Range("D25").Value = "WINYBEL SANDRY DE LA HOZ MUÑOZ"
$ws.Range("E25").Value = "2104"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 908526

$ws.Range("C26").Value = "1051890181"
$ws.Range("D26").Value = "WINYBEL SANDRY DE LA HOZ MUÑOZ"
$ws.Range("E26").Value = "2103"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 908526

# HERNANDO MEDINA SANCHEZ (CC 91427944) - rows 27-33
$ws.Range("C27").Value = "91427944"
$ws.Range("D27").Value = "HERNANDO MEDINA SANCHEZ"
$ws.Range("E27").Value = "2109"
$ws.Range("F27").Value = 53334
$ws.Range("G27").Value = 2000000

$ws.Range("C28").Value = "91427944"
$ws.Range("D28").Value = "HERNANDO MEDINA SANCHEZ"
$ws.Range("E28").Value = "2108"
$ws.Range("F28").Value = 80000
$ws.Range("G28").Value = 2000000

$ws.Range("C29").Value = "91427944"
$ws.Range("D29").Value = "HERNANDO MEDINA SANCHEZ"
$ws.Range("E29").Value = "2107"
$ws.Range("F29").Value = 80000
$ws.Range("G29").Value = 2000000

$ws.Range("C30").Value = "91427944"
$ws.Range("D30").Value = "HERNANDO MEDINA SANCHEZ"
$ws.Range("E30").Value = "2106"
$ws.Range("F30").Value = 80000
$ws.Range("G30").Value = 2000000

$ws.Range("C31").Value = "91427944"
$ws.Range("D31").Value = "HERNANDO MEDINA SANCHEZ"
$ws.Range("E31").Value = "2105"
$ws.Range("F31").Value = 80000
$ws.Range("G31").Value = 2000000

$ws.Range("C32").Value = "91427944"
$ws.Range("D32").Value = "HERNANDO MEDINA SANCHEZ"
$ws.Range("E32").Value = "2104"
$ws.Range("F32").Value = 80000
$ws.Range("G32").Value = 2000000

$ws.Range("C33").Value = "91427944"
$ws.Range("D33").Value = "HERNANDO MEDINA SANCHEZ"
$ws.Range("E33").Value = "2103"
$ws.Range("F33").Value = 80000
$ws.Range("G33").Value = 2000000
